{"js": "// Replace review copy text as described by the diff.\n// Each entry: exact source text -> exact replacement text.\nconst replacements = [\n  [\n    \"Play Mystic Mirror Free Slot | Beautiful Fairy Tale Theme\",\n    \"Play Mystic Mirror Free - Exciting Fairy-Tale Slot Game\",\n  ],\n  [\n    \"Sixth reel makes it easier to form larger and better-paying combinations\",\n    \"Unique gameplay with the addition of a sixth reel\",\n  ],\n  [\n    \"Special symbols include Wild/Scatter and Magic Symbols that can expand on the reels\",\n    \"Special symbols and bonus features\",\n  ],\n  [\n    \"Free Spins feature with random symbols that can cover the entire screen\",\n    \"10 jackpot prizes during Free Spins\",\n  ],\n  [\n    \"Does not feature as many exciting bonus features as other titles offered by Red Rake Gaming\",\n    \"Not as many exciting features compared to other Red Rake Gaming titles\",\n  ],\n  [\n    \"Betting system is rather unusual and can be confusing for some players\",\n    \"May not be the most successful game by Red Rake\",\n  ],\n  [\n    \"Explore Mystic Mirror, an online slot game with a fairy tale theme and special features like Free Spins and Magic Symbols. Play for free today.\",\n    \"Read our review of Mystic Mirror, an online slot game with a magical fairy-tale theme. Play for free and enjoy unique gameplay and special features.\",\n  ],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the review copy rewrite described by the diff using Find/Replace.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"Play Mystic Mirror Free Slot | Beautiful Fairy Tale Theme\"; Replace = \"Play Mystic Mirror Free - Exciting Fairy-Tale Slot Game\" },\n    @{ Find = \"Sixth reel makes it easier to form larger and better-paying combinations\"; Replace = \"Unique gameplay with the addition of a sixth reel\" },\n    @{ Find = \"Special symbols include Wild/Scatter and Magic Symbols that can expand on the reels\"; Replace = \"Special symbols and bonus features\" },\n    @{ Find = \"Free Spins feature with random symbols that can cover the entire screen\"; Replace = \"10 jackpot prizes during Free Spins\" },\n    @{ Find = \"Does not feature as many exciting bonus features as other titles offered by Red Rake Gaming\"; Replace = \"Not as many exciting features compared to other Red Rake Gaming titles\" },\n    @{ Find = \"Betting system is rather unusual and can be confusing for some players\"; Replace = \"May not be the most successful game by Red Rake\" },\n    @{ Find = \"Explore Mystic Mirror, an online slot game with a fairy tale theme and special features like Free Spins and Magic Symbols. Play for free today.\"; Replace = \"Read our review of Mystic Mirror, an online slot game with a magical fairy-tale theme. Play for free and enjoy unique gameplay and special features.\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    $find.Execute(\n        $find.Text,\n        $false,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $find.Replacement.Text,\n        2\n    )\n}\n"}
